$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 494, pushing the existing rows
# 494..572 down to 496..574 (preserving their data/styles).
$ws.Rows.Item(494).Resize(2).Insert()

# New row 494 - weekly update row (Primera, $/caja 36 atados)
$ws.Cells.Item(494,1).Value = 6
$ws.Cells.Item(494,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(494,3).Value = "Metropolitana"
$ws.Cells.Item(494,4).Value = 44505
$ws.Cells.Item(494,5).Value = 13
$ws.Cells.Item(494,6).Value = 100112040
$ws.Cells.Item(494,7).Value = "Cilantro"
$ws.Cells.Item(494,8).Value = "Sin especificar"
$ws.Cells.Item(494,9).Value = "Primera"
$ws.Cells.Item(494,10).Value = 720
$ws.Cells.Item(494,11).Value = 3500
$ws.Cells.Item(494,12).Value = 4000
$ws.Cells.Item(494,13).Value = 3729
$ws.Cells.Item(494,14).Value = "$/caja 36 atados"
$ws.Cells.Item(494,15).Value = "Región Metropolitana"
$ws.Cells.Item(494,16).Value = 104
$ws.Cells.Item(494,17).Value = 36
$ws.Cells.Item(494,18).Value = "Hortaliza"

# New row 495 - weekly update row (Primera, $/docena de atados)
$ws.Cells.Item(495,1).Value = 6
$ws.Cells.Item(495,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(495,3).Value = "Metropolitana"
$ws.Cells.Item(495,4).Value = 44505
$ws.Cells.Item(495,5).Value = 13
$ws.Cells.Item(495,6).Value = 100112040
$ws.Cells.Item(495,7).Value = "Cilantro"
$ws.Cells.Item(495,8).Value = "Sin especificar"
$ws.Cells.Item(495,9).Value = "Primera"
$ws.Cells.Item(495,10).Value = 430
$ws.Cells.Item(495,11).Value = 6500
$ws.Cells.Item(495,12).Value = 7000
$ws.Cells.Item(495,13).Value = 6674
$ws.Cells.Item(495,14).Value = "$/docena de atados"
$ws.Cells.Item(495,15).Value = "Región Metropolitana"
$ws.Cells.Item(495,16).Value = 2225
$ws.Cells.Item(495,17).Value = 3
$ws.Cells.Item(495,18).Value = "Hortaliza"
